# JobMaterial.xlsx ("more estimates") edit.
#
# The "Banner" element rows (2-5, originally Quantity 0.08) and the
# "Flat  2p" element rows (6-9, originally Quantity 0.07) swap places:
# rows 2-5 become "Flat  2p" @ 0.07 and rows 6-9 become "Banner" @ 0.08.
# In addition, the Wide Format UV color (column D) used for the
# now-"Banner" rows 6-9 is rotated by one position
# (Yellow -> Black -> Cyan -> Magenta -> Yellow).
#
# Columns E2:E9 ("Quantity") are stored as TEXT in the workbook (e.g. the
# literal string "0.07", not the number 0.07) - format the cells as Text
# first so the COM layer doesn't auto-coerce the numeric-looking string
# into a number, then restore the original font so appearance is
# unaffected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-QuantityText($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Font.Name = "Calibri"
    $rng.Font.Size = 14
    $rng.Font.Color = 255
    $rng.Value = $val
}

# Rows 2-5: "Banner" @ 0.08  ->  "Flat  2p" @ 0.07 (colors unchanged)
$ws.Range("B2").Value = "Flat  2p"
Set-QuantityText "E2" "0.07"

$ws.Range("B3").Value = "Flat  2p"
Set-QuantityText "E3" "0.07"

$ws.Range("B4").Value = "Flat  2p"
Set-QuantityText "E4" "0.07"

$ws.Range("B5").Value = "Flat  2p"
Set-QuantityText "E5" "0.07"

# Rows 6-9: "Flat  2p" @ 0.07  ->  "Banner" @ 0.08, with the color
# (Material, column D) rotated by one: Yellow->Black, Black->Cyan,
# Cyan->Magenta, Magenta->Yellow
$ws.Range("B6").Value = "Banner"
$ws.Range("D6").Value = "Black - Wide Format UV - "
Set-QuantityText "E6" "0.08"

$ws.Range("B7").Value = "Banner"
$ws.Range("D7").Value = "Cyan - Wide Format UV - "
Set-QuantityText "E7" "0.08"

$ws.Range("B8").Value = "Banner"
$ws.Range("D8").Value = "Magenta - Wide Format UV - "
Set-QuantityText "E8" "0.08"

$ws.Range("B9").Value = "Banner"
$ws.Range("D9").Value = "Yellow - Wide Format UV - "
Set-QuantityText "E9" "0.08"

# Rows 10-11: Element values stay the same content ("Flat  2p" / "Banner")
$ws.Range("B10").Value = "Flat  2p"
$ws.Range("B11").Value = "Banner"
